$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 111: 2020-02-27
$ws.Range("A111").Value = 1582761600
$ws.Range("B111").NumberFormat = "@"
$ws.Range("B111").Value = "2020-02-27"
$ws.Range("B111").ClearFormats()
$ws.Range("C111").NumberFormat = "@"
$ws.Range("C111").Value = "03027"
$ws.Range("C111").ClearFormats()
$ws.Range("D111").Value = "MMIS"
$ws.Range("E111").Value = 0.165
$ws.Range("F111").Value = 0.165
$ws.Range("G111").Value = 0.165
$ws.Range("H111").Value = 0.165
$ws.Range("I111").NumberFormat = "@"
$ws.Range("I111").Value = "-"
$ws.Range("I111").ClearFormats()

# Row 112: 2020-02-28
$ws.Range("A112").Value = 1582848000
$ws.Range("B112").NumberFormat = "@"
$ws.Range("B112").Value = "2020-02-28"
$ws.Range("B112").ClearFormats()
$ws.Range("C112").NumberFormat = "@"
$ws.Range("C112").Value = "03027"
$ws.Range("C112").ClearFormats()
$ws.Range("D112").Value = "MMIS"
$ws.Range("E112").Value = 0.165
$ws.Range("F112").Value = 0.165
$ws.Range("G112").Value = 0.165
$ws.Range("H112").Value = 0.165
$ws.Range("I112").NumberFormat = "@"
$ws.Range("I112").Value = "-"
$ws.Range("I112").ClearFormats()
